$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Gal"
$ws.Range("C2").Value = "Gpr151"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.48067
$ws.Range("H2").Value = 1.44201
$ws.Range("I2").Value = 0.1949338371837906
$ws.Range("J2").Value = 0.1949338371837907
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3933823333333333
$ws.Range("N2").Value = 1.180147
$ws.Range("O2").Value = 0.3432826356401089
$ws.Range("P2").Value = 0.3432826356401089
$ws.Range("Q2").Value = 0.1890870861633333
$ws.Range("R2").Value = 1.70178377547
$ws.Range("S2").Value = 0.06691740140389152
$ws.Range("T2").Value = 0.06691740140389153
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Gal"
$ws.Range("C3").Value = "Gpr151"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.48067
$ws.Range("H3").Value = 1.44201
$ws.Range("I3").Value = 0.1949338371837906
$ws.Range("J3").Value = 0.1949338371837907
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.4983106666666666
$ws.Range("N3").Value = 1.494932
$ws.Range("O3").Value = 0.4348476902137948
$ws.Range("P3").Value = 0.4348476902137948
$ws.Range("Q3").Value = 0.2395229881466666
$ws.Range("R3").Value = 2.15570689332
$ws.Range("S3").Value = 0.08476652884388332
$ws.Range("T3").Value = 0.08476652884388333
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gal"
$ws.Range("C4").Value = "Gpr151"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.48067
$ws.Range("H4").Value = 1.44201
$ws.Range("I4").Value = 0.1949338371837906
$ws.Range("J4").Value = 0.1949338371837907
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.25425
$ws.Range("N4").Value = 0.7627499999999999
$ws.Range("O4").Value = 0.2218696741460963
$ws.Range("P4").Value = 0.2218696741460963
$ws.Range("Q4").Value = 0.1222103475
$ws.Range("R4").Value = 1.0998931275
$ws.Range("S4").Value = 0.04324990693601581
$ws.Range("T4").Value = 0.04324990693601582
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Gal"
$ws.Range("C5").Value = "Gpr151"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.985141
$ws.Range("H5").Value = 5.955423
$ws.Range("I5").Value = 0.8050661628162092
$ws.Range("J5").Value = 0.8050661628162094
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3933823333333333
$ws.Range("N5").Value = 1.180147
$ws.Range("O5").Value = 0.3432826356401089
$ws.Range("P5").Value = 0.3432826356401089
$ws.Range("Q5").Value = 0.7809193985756665
$ws.Range("R5").Value = 7.028274587180999
$ws.Range("S5").Value = 0.2763652342362173
$ws.Range("T5").Value = 0.2763652342362174
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Gal"
$ws.Range("C6").Value = "Gpr151"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.985141
$ws.Range("H6").Value = 5.955423
$ws.Range("I6").Value = 0.8050661628162092
$ws.Range("J6").Value = 0.8050661628162094
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.4983106666666666
$ws.Range("N6").Value = 1.494932
$ws.Range("O6").Value = 0.4348476902137948
$ws.Range("P6").Value = 0.4348476902137948
$ws.Range("Q6").Value = 0.9892169351373331
$ws.Range("R6").Value = 8.902952416235999
$ws.Range("S6").Value = 0.3500811613699115
$ws.Range("T6").Value = 0.3500811613699115
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Gal"
$ws.Range("C7").Value = "Gpr151"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.985141
$ws.Range("H7").Value = 5.955423
$ws.Range("I7").Value = 0.8050661628162092
$ws.Range("J7").Value = 0.8050661628162094
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.25425
$ws.Range("N7").Value = 0.7627499999999999
$ws.Range("O7").Value = 0.2218696741460963
$ws.Range("P7").Value = 0.2218696741460963
$ws.Range("Q7").Value = 0.5047220992499999
$ws.Range("R7").Value = 4.542498893249999
$ws.Range("S7").Value = 0.1786197672100804
$ws.Range("T7").Value = 0.1786197672100804

$ws.Rows("8:10").Delete()
